$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet has two "section header" rows that carry only a label in
# column A (no data in B:I): row 8 "grandes regiões e unidades da
# federação" and row 5 "situação do domicílio". Both are being removed,
# which shifts every row below each of them up by one.
#
# Delete the lower one (row 8) first so the row number of the upper one
# (row 5) is not affected by the first deletion.
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(5).Delete()

# Row 2 held a spurious "unnamed: 1_level_1" label over the "total"
# column and was missing the real "total" label; fix the row-2 header so
# B2:I2 read total / 10 a 15 anos / 16 a 24 anos / ... / 70 anos ou mais.
$ws.Range("B2").Value = "total"
$ws.Range("C2").Value = "10 a 15 anos"
$ws.Range("D2").Value = "16 a 24 anos"
$ws.Range("E2").Value = "25 a 34 anos"
$ws.Range("F2").Value = "35 a 49 anos"
$ws.Range("G2").Value = "50 a 59 anos"
$ws.Range("H2").Value = "60 a 69 anos"
$ws.Range("I2").Value = "70 anos ou mais"
